# Update gh-pages to output generated at 456a3b4
# Refresh the "想去人数" (interested-attendee count) figures pulled from
# bilibili for a handful of events. The same events appear both on the
# "展览" (Exhibitions) sheet and on the aggregated "全部类型" (All types)
# sheet, so each figure has to be updated in both places. "全部类型" has
# one extra (non-exhibition) row inserted above, so the matching rows are
# offset by one there.

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions) sheet
$wsExhibitions = $wb.Worksheets.Item("展览")
$wsExhibitions.Range("F2").Value = 6039    # 南昌·Sunflower Garden动漫游戏展: 5531 -> 6039
$wsExhibitions.Range("F6").Value = 50      # 鹰潭·MZD动漫游戏嘉年华: 49 -> 50
$wsExhibitions.Range("F9").Value = 32      # 赣州·十万伏特-第七届青年文化综合展览会: 31 -> 32
$wsExhibitions.Range("F17").Value = 169    # 景德镇·第十六届瓷都ACG内场—花玲&宴宁: 168 -> 169
$wsExhibitions.Range("F18").Value = 1668   # 江西·JMG（江西广电）第二届UP动漫游戏博览会: 1666 -> 1668

# 全部类型 (All types) sheet - same events, rows shifted by +1 from row 9 on
$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F2").Value = 6039       # 南昌·Sunflower Garden动漫游戏展: 5531 -> 6039
$wsAllTypes.Range("F6").Value = 50         # 鹰潭·MZD动漫游戏嘉年华: 49 -> 50
$wsAllTypes.Range("F10").Value = 32        # 赣州·十万伏特-第七届青年文化综合展览会: 31 -> 32
$wsAllTypes.Range("F18").Value = 169       # 景德镇·第十六届瓷都ACG内场—花玲&宴宁: 168 -> 169
$wsAllTypes.Range("F19").Value = 1668      # 江西·JMG（江西广电）第二届UP动漫游戏博览会: 1666 -> 1668
